$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.203.56'
$ws.Range("E2").Value = '  +6.22%  '
$ws.Range("D3").Value = '2.437.96'
$ws.Range("E3").Value = '  +6.11%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '565.10'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +4.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.69'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +11.64%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  +3.75%  '
$ws.Range("D9").Value = '2.435.14'
$ws.Range("E9").Value = '  +6.08%  '
$ws.Range("E10").Value = '  +4.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.75'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +5.05%  '
$ws.Range("E12").Value = '  +1.34%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.352'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +7.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.38'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +14.88%  '
$ws.Range("D15").Value = '2.872.79'
$ws.Range("E15").Value = '  +6.15%  '
$ws.Range("D16").Value = '63.092.17'
$ws.Range("E16").Value = '  +6.29%  '
$ws.Range("E17").Value = '  +9.47%  '
$ws.Range("D18").Value = '2.438.58'
$ws.Range("E18").Value = '  +6.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.21'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +8.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '340.16'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +10.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.28'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +6.78%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.73'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +4.28%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.26'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +4.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.174'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +3.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("E27").Value = '  +14.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.15'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +6.71%  '
$ws.Range("E29").Value = '  +13.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.69'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +16.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.82'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +8.08%  '
$ws.Range("D32").Value = '0.0₃0787'
$ws.Range("E32").Value = '  +11.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '174.53'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.83%  '
$ws.Range("E34").Value = '  +12.61%  '
$ws.Range("E35").Value = '  +6.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.71'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +6.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '371.66'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +20.52%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.48'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +13.34%  '
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.71'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +14.72%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '40.34'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +7.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '149.39'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +10.76%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.69'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +9.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.73'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +13.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.595'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +5.68%  '
$ws.Range("E47").Value = '  +3.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0520'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +7.03%  '
$ws.Range("E49").Value = '  +7.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.90'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +8.28%  '
$ws.Range("D51").Value = '0.0₆0223'
$ws.Range("E51").Value = '  +2.70%  '
